$wb = $excel.ActiveWorkbook

# --- Brand sheet: add new row (A6, B6) ---
$wsBrand = $wb.Worksheets.Item("Brand")
$wsBrand.Cells.Item(6, 1).Value = 5
$wsBrand.Cells.Item(6, 2).Value = "yuyuy"

# --- Category sheet: update existing cell B2 ---
$wsCategory = $wb.Worksheets.Item("Category")
$wsCategory.Cells.Item(2, 2).Value = "sff"

# --- Product sheet: add new row (A3:G3) ---
$wsProduct = $wb.Worksheets.Item("Product")
$wsProduct.Cells.Item(3, 1).Value = 2
$wsProduct.Cells.Item(3, 2).Value = "df"
$wsProduct.Cells.Item(3, 3).Value = "Categoría 1"
$wsProduct.Cells.Item(3, 4).Value = "Marca B"
$wsProduct.Cells.Item(3, 5).Value = 11
$wsProduct.Cells.Item(3, 6).Value = 1
# Leading apostrophe forces this numeric-looking value to be stored as text
# (matches the source data, which keeps "Location" values like "1" as strings).
$wsProduct.Cells.Item(3, 7).Value = "'1"
